{"js": "// Update the run metadata (Start time, End time, Duration) in the body text.\nconst replacements = [\n  { find: \"Start time: 2017-12-27 18:33:20\", replace: \"Start time: 2018-01-31 12:36:43\" },\n  { find: \"End time: 2017-12-27 18:33:26\", replace: \"End time: 2018-01-31 12:36:51\" },\n  { find: \"Duration: 6.57 secs\", replace: \"Duration: 8.53 secs\" }\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\nReplace-Text \"Start time: 2017-12-27 18:33:20\" \"Start time: 2018-01-31 12:36:43\"\nReplace-Text \"End time: 2017-12-27 18:33:26\" \"End time: 2018-01-31 12:36:51\"\nReplace-Text \"Duration: 6.57 secs\" \"Duration: 8.53 secs\"\n"}
